# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (F) counts and "最低票价" (G) values on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 3146

    if ($name -eq "展览") {
        $ws.Range("F4").Value = 114
        $ws.Range("F5").Value = 6852
        $ws.Range("F6").Value = 1921
        $ws.Range("G6").Value = 65
        $ws.Range("F7").Value = 12
        $ws.Range("F8").Value = 68
        $ws.Range("F14").Value = 166
        $ws.Range("F15").Value = 31
    }
    else {
        $ws.Range("F5").Value = 114
        $ws.Range("F6").Value = 6852
        $ws.Range("F7").Value = 1921
        $ws.Range("G7").Value = 65
        $ws.Range("F8").Value = 12
        $ws.Range("F9").Value = 68
        $ws.Range("F15").Value = 166
        $ws.Range("F16").Value = 31
    }
}
